# 349. Intersection of Two Array and GFG. Union of Two sorted array
#
# Row 25 ("Searching an element in a sorted array") previously carried a
# one-off "heading-ish" look (taller row, bold 13.5pt question cell).
# The edit normalizes row 25 back to the plain look used by every other
# data row, then appends two new rows for the two new questions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize row 25's formatting (drop the bold/oversized font and the
#     explicit row height) so it matches the other plain rows. ---
$ws.Rows.Item(25).AutoFit() | Out-Null
$ws.Range("C25").Copy() | Out-Null
$ws.Range("B25").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Row 26: GFG / Union of Two Sorted Arrays / Java / same date as above ---
$ws.Range("A25").Copy() | Out-Null
$ws.Range("A26").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("B25").Copy() | Out-Null
$ws.Range("B26").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("C25").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("D25").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A26").Value = "GFG"
$ws.Range("B26").Value = "Union of Two Sorted Arrays"
$ws.Range("C26").Value = "Java"
$ws.Range("D26").Value = $ws.Range("D25").Value2

# --- Row 27: 349 / Intersection of Two Arrays / Java / same date as above ---
$ws.Range("A25").Copy() | Out-Null
$ws.Range("A27").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("B25").Copy() | Out-Null
$ws.Range("B27").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("C25").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("D25").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A27").Value = 349
$ws.Range("B27").Value = "Intersection of Two Arrays"
$ws.Range("C27").Value = "Java"
$ws.Range("D27").Value = $ws.Range("D25").Value2

$ws.Application.CutCopyMode = $false

# --- Match the author's resulting view/selection state ---
$ws.Range("M14").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
